$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.453.76'
$ws.Range("E2").Value = '  +0.39%  '
$ws.Range("D3").Value = '1.828.78'
$ws.Range("E3").Value = '  +0.16%  '
$ws.Range("E4").Value = '  +0.56%  '
$ws.Range("D5").Value = "'314.84"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.84%  '
$ws.Range("E6").Value = '  +0.48%  '
$ws.Range("D7").Value = "'0.5122"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -4.06%  '
$ws.Range("D8").Value = "'0.3926"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -1.66%  '
$ws.Range("D9").Value = "'0.07668"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +1.42%  '
$ws.Range("D10").Value = "'41.76"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.12%  '
$ws.Range("E11").Value = '  +0.95%  '
$ws.Range("D12").Value = "'21.09"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +2.14%  '
$ws.Range("D13").Value = "'6.313"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.64%  '
$ws.Range("D14").Value = "'1.002"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.53%  '
$ws.Range("D15").Value = "'7.539"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -1.08%  '
$ws.Range("D16").Value = '1.828.85'
$ws.Range("E16").Value = '  +0.96%  '
$ws.Range("D17").Value = "'93.72"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +4.65%  '
$ws.Range("D18").Value = "'0.00001102"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +3.31%  '
$ws.Range("D19").Value = "'0.06706"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +1.80%  '
$ws.Range("D20").Value = "'17.69"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.61%  '
$ws.Range("E21").Value = '  +0.37%  '
$ws.Range("D22").Value = "'6.148"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +1.86%  '
$ws.Range("D23").Value = '28.488.54'
$ws.Range("E23").Value = '  +0.46%  '
$ws.Range("E24").Value = '  -0.23%  '
$ws.Range("D25").Value = "'2.257"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +8.08%  '
$ws.Range("D26").Value = "'20.80"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +1.34%  '
$ws.Range("D27").Value = "'156.92"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.47%  '
$ws.Range("D28").Value = '2.038.99'
$ws.Range("E28").Value = '  +0.91%  '
$ws.Range("D29").Value = "'2.400"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.32%  '
$ws.Range("D30").Value = "'124.54"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.56%  '
$ws.Range("D31").Value = "'1.117"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.70%  '
$ws.Range("D32").Value = "'0.1085"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -1.09%  '
$ws.Range("D33").Value = "'5.671"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +1.53%  '
$ws.Range("D34").Value = "'3.664"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.15%  '
$ws.Range("D35").Value = "'0.07035"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -5.26%  '
$ws.Range("E36").Value = '  -0.69%  '
$ws.Range("D37").Value = "'8.964"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +3.66%  '
$ws.Range("D38").Value = "'0.02323"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.90%  '
$ws.Range("D39").Value = "'5.154"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -1.24%  '
$ws.Range("D40").Value = "'0.6277"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.90%  '
$ws.Range("D41").Value = "'11.22"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.98%  '
$ws.Range("D42").Value = "'1.179"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -1.04%  '
$ws.Range("D43").Value = "'1.000"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.21%  '
$ws.Range("D44").Value = "'1.390"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.91%  '
$ws.Range("D45").Value = "'13.45"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.29%  '
$ws.Range("D46").Value = "'0.5899"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +1.81%  '
$ws.Range("D47").Value = "'3.715"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.76%  '
$ws.Range("D48").Value = "'124.75"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.34%  '
$ws.Range("D49").Value = "'1.980"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +1.43%  '
$ws.Range("E50").Value = '  +0.85%  '
$ws.Range("E51").Value = '  +0.87%  '
